$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.736.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.085.37'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.66%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '540.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.46%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.078.39'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.496'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.46%  '

$ws.Range("E10").Value = '  -2.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.25'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.89%  '

$ws.Range("E12").Value = '  +0.52%  '

$ws.Range("E13").Value = '  +4.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.584.85'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.37%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.712.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.12%  '

$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.088.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("E19").Value = '  +0.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '488.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.702'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("E27").Value = '  -0.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("E30").Value = '  -0.33%  '

$ws.Range("E31").Value = '  -2.83%  '

$ws.Range("E32").Value = '  +0.15%  '

$ws.Range("E33").Value = '  -4.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.22'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.47'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.81%  '

$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.12%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '494.54'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.274.11'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0400'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.81%  '

$ws.Range("E40").Value = '  +0.80%  '

$ws.Range("E41").Value = '  -1.23%  '

$ws.Range("E42").Value = '  +0.68%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.22%  '

$ws.Range("E44").Value = '  +1.13%  '

$ws.Range("E46").Value = '  +0.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0539'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.06%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.06%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.44%  '

$ws.Range("E50").Value = '  +2.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.37'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.17%  '
